# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple Price (D) / Volume(1h) (E) updates -------------------------------
# Row -> Price, Volume
$updates = @(
    @{ Row = 2;  D = "72.070.55";  E = "  +0.42%  " }
    @{ Row = 3;  D = "4.040.16";   E = "  -0.11%  " }
    @{ Row = 4;  E = "  +0.02%  " }
    @{ Row = 5;  D = "538.04";     E = "  +0.86%  " }
    @{ Row = 6;  D = "149.37";     E = "  -3.28%  " }
    @{ Row = 7;  D = "4.036.76";   E = "  +0.00%  " }
    @{ Row = 8;  D = "0.697";      E = "  +0.36%  " }
    @{ Row = 9;  D = "1.00";       E = "  -0.05%  " }
    @{ Row = 10; E = "  -1.09%  " }
    @{ Row = 11; E = "  -2.33%  " }
    @{ Row = 12; D = "53.55";      E = "  +8.16%  " }
    @{ Row = 13; E = "  -1.19%  " }
    @{ Row = 14; D = "10.91";      E = "  -0.92%  " }
    @{ Row = 15; D = "4.685.58";   E = "  -0.05%  " }
    @{ Row = 16; D = "4.031.17";   E = "  -0.23%  " }
    @{ Row = 17; D = "14.26";      E = "  -1.27%  " }
    @{ Row = 18; E = "  -1.60%  " }
    @{ Row = 19; E = "  -1.78%  " }
    @{ Row = 20; E = "  -1.30%  " }
    @{ Row = 21; D = "72.104.37";  E = "  +0.42%  " }
    @{ Row = 22; D = "436.96";     E = "  -0.10%  " }
    @{ Row = 23; D = "98.22";      E = "  -1.58%  " }
    @{ Row = 24; D = "3.52";       E = "  -5.34%  " }
    @{ Row = 25; E = "  +0.42%  " }
    @{ Row = 26; D = "14.66";      E = "  -1.40%  " }
    @{ Row = 27; D = "4.37";       E = "  +23.31%  " }
    @{ Row = 28; D = "11.28";      E = "  -1.30%  " }
    @{ Row = 29; D = "10.76";      E = "  -1.64%  " }
    @{ Row = 30; D = "5.96";       E = "  +2.07%  " }
    @{ Row = 31; D = "37.11";      E = "  -0.66%  " }
    @{ Row = 32; D = "8.34";       E = "  +23.45%  " }
    @{ Row = 33; E = "  +2.32%  " }
    @{ Row = 34; D = "50.28";      E = "  +16.53%  " }
    @{ Row = 35; D = "13.60";      E = "  -0.70%  " }
    @{ Row = 36; D = "682.33";     E = "  +1.30%  " }
    @{ Row = 37; D = "67.30";      E = "  +0.85%  " }
    @{ Row = 38; D = "0.462";      E = "  +6.30%  " }
    @{ Row = 39; D = "0.0₃0869";   E = "  +1.15%  " }
    @{ Row = 42; E = "  -0.82%  " }
    @{ Row = 43; D = "11.20";      E = "  +17.44%  " }
    @{ Row = 44; E = "  +0.03%  " }
    @{ Row = 45; E = "  -1.02%  " }
    @{ Row = 46; D = "1.00";       E = "  +0.20%  " }
    @{ Row = 47; E = "  -0.87%  " }
    @{ Row = 48; E = "  -2.68%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("D")) {
        $ws.Range("D$r").NumberFormat = "@"
        $ws.Range("D$r").Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E$r").Value = $u.E
    }
}

# --- Rows whose Coin/Link/Price/Volume got re-ordered/replaced ---------------
# Row 40 <- WEMIXToken (was on row 41)
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.44"
$ws.Range("E40").Value = "  +7.84%  "

# Row 41 <- Kaspa (was on row 40)
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.149"
$ws.Range("E41").Value = "  -5.83%  "

# Row 49 <- Stacks (was on row 50)
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.11"
$ws.Range("E49").Value = "  +1.39%  "

# Row 50 <- ApeXProtocol (was on row 49)
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.32"
$ws.Range("E50").Value = "  -1.68%  "

# Row 51 <- FLOKI (was Maker)
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000281"
$ws.Range("E51").Value = "  +2.68%  "
